$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so numeric-looking
# values such as "597.96" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.352.55"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.649.78"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D5").Value = "597.96"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "159.64"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "28.09"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "3.130.35"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "68.334.87"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.626.04"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "11.42"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "364.57"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.43"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.33"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "75.11"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "9.72"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").Value = "2.784.18"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.03"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "558.57"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "1.87"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "1.57"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "19.90"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("D38").Value = "159.63"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").Value = "5.36"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "0.0₆0336"
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "17.80"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D46").Value = "158.01"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "3.77"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "22.27"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  +0.02%  "
